$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$row1 = 184
$row2 = 185

$ws.Cells.Item($row1, 1).Value = 4
$ws.Cells.Item($row1, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item($row1, 3).Value = "Los Lagos"
$ws.Cells.Item($row1, 4).Value = 44656
$ws.Cells.Item($row1, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item($row1, 5).Value = 10
$ws.Cells.Item($row1, 6).Value = "Fruta"
$ws.Cells.Item($row1, 7).Value = 100103
$ws.Cells.Item($row1, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item($row1, 9).Value = 100103002
$ws.Cells.Item($row1, 10).Value = "Ciruela"
$ws.Cells.Item($row1, 11).Value = "Angeleno"
$ws.Cells.Item($row1, 12).Value = "Primera"
$ws.Cells.Item($row1, 13).Value = 600
$ws.Cells.Item($row1, 14).Value = 14000
$ws.Cells.Item($row1, 15).Value = 15000
$ws.Cells.Item($row1, 16).Value = 14500
$ws.Cells.Item($row1, 17).Value = "$/caja 15 kilos granel"
$ws.Cells.Item($row1, 18).Value = "Provincia de Curicó"
$ws.Cells.Item($row1, 19).Value = 967
$ws.Cells.Item($row1, 20).Value = 15

$ws.Cells.Item($row2, 1).Value = 4
$ws.Cells.Item($row2, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item($row2, 3).Value = "Los Lagos"
$ws.Cells.Item($row2, 4).Value = 44656
$ws.Cells.Item($row2, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item($row2, 5).Value = 10
$ws.Cells.Item($row2, 6).Value = "Fruta"
$ws.Cells.Item($row2, 7).Value = 100103
$ws.Cells.Item($row2, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item($row2, 9).Value = 100103002
$ws.Cells.Item($row2, 10).Value = "Ciruela"
$ws.Cells.Item($row2, 11).Value = "Angeleno"
$ws.Cells.Item($row2, 12).Value = "Segunda"
$ws.Cells.Item($row2, 13).Value = 300
$ws.Cells.Item($row2, 14).Value = 13000
$ws.Cells.Item($row2, 15).Value = 13000
$ws.Cells.Item($row2, 16).Value = 13000
$ws.Cells.Item($row2, 17).Value = "$/caja 15 kilos granel"
$ws.Cells.Item($row2, 18).Value = "Provincia de Curicó"
$ws.Cells.Item($row2, 19).Value = 867
$ws.Cells.Item($row2, 20).Value = 15
